# "add c++ json parsing and additions to sbc list and formations"
#
# This edits the SBC (Squad Building Challenge) list on Sheet1:
#  - the two placeholder rows (4 and 5), which were just duplicates of the
#    row-2 "league and nation hybrid" entry, are removed
#  - row 3 (previously another duplicate placeholder) is replaced with a
#    real new SBC entry: "Icons" / "Petit" / "84-Rated Squad"
#  - the B:C columns are widened so the new text fits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing placeholder rows (old rows 4 & 5).
$ws.Range("A4:A5").EntireRow.Delete() | Out-Null

# Fill in row 3 with the new SBC entry.
# (name column, "C", is set first so its new string lands in the shared
# string table ahead of the other two new strings, matching source order.)
$ws.Range("C3").Value = "84-Rated Squad"
$ws.Range("A3").Value = "Icons"
$ws.Range("B3").Value = "Petit"
$ws.Range("D3").Value = 41212
$ws.Range("P3").Value = 75
$ws.Range("Q3").Value = 84
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 11
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 11
$ws.Range("V3").Value = 11
$ws.Range("W3").Value = 5

# Widen the "sub category" / "name" columns so the new values fit.
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(3).ColumnWidth = 14

# Move the active selection up to A2.
$ws.Range("A2").Select() | Out-Null
